$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row counts (column K) to reflect additional simulation samples ---
$kUpdates = @{
    13 = 37911
    14 = 4278
    15 = 1104
    16 = 1012
    17 = 4620
    18 = 1310
    19 = 2939
    20 = 509
    21 = 216
    22 = 4329
    23 = 126
    27 = 182
    29 = 4527
    30 = 636
    31 = 172
    32 = 8360
    33 = 1368
    34 = 391
    36 = 551
    37 = 368
    38 = 657
    39 = 1213
    40 = 403
    41 = 4459
    42 = 745
    43 = 638
    44 = 817
    45 = 715
    48 = 664
    50 = 1618
    51 = 382
    52 = 668
    53 = 96
    57 = 105
    58 = 1153
    59 = 1298
    60 = 1347
    61 = 396
    62 = 527
    63 = 255
    65 = 366
    66 = 1175
    67 = 203
    68 = 203
    69 = 812
    70 = 540
    71 = 388
    72 = 589
    73 = 1491
    79 = 1335
    80 = 329
    81 = 321
    82 = 272
    83 = 659
    84 = 4352
    86 = 1360
    87 = 349
    91 = 186
    92 = 413
    93 = 334
    95 = 318
    96 = 702
    97 = 137
    98 = 667
    99 = 768
    100 = 556
    102 = 1303
    103 = 329
    104 = 122
    105 = 668
    106 = 594
    107 = 292
    108 = 415
    109 = 351
    110 = 107
    111 = 107
    112 = 146
    113 = 831
    114 = 1327
    116 = 451
    117 = 123
    118 = 895
    119 = 354
    121 = 271
    122 = 676
    123 = 403
    124 = 706
    127 = 572
    128 = 151
    129 = 370
    130 = 793
    131 = 398
    132 = 219
    133 = 331
    134 = 142
    136 = 383
    137 = 100
    138 = 98
    140 = 392
    141 = 144
    142 = 504
    143 = 395
    145 = 387
    148 = 365
    149 = 626
    150 = 207
    155 = 189
    156 = 1426
    157 = 242
    158 = 539
    159 = 166
    160 = 187
    161 = 374
    162 = 655
    163 = 275
    165 = 778
    166 = 274
    168 = 136
    169 = 345
    171 = 639
    173 = 144
    174 = 128
    175 = 207
    176 = 75
    177 = 189
    178 = 179
    180 = 1305
    181 = 305
    182 = 168
    183 = 119
    184 = 135
    186 = 127
    187 = 186
    188 = 84
    189 = 149
    190 = 149
    191 = 215
    192 = 288
    195 = 205
    196 = 97
    197 = 104
    198 = 103
    200 = 219
    201 = 179
    202 = 215
    203 = 334
    204 = 324
    205 = 136
    207 = 94
    208 = 86
    209 = 121
    210 = 69
    211 = 336
    212 = 54
    213 = 110
    214 = 493
    215 = 113
    216 = 125
    217 = 116
    230 = 51
    259 = -15
    407 = 3
    692 = 4
}
foreach ($row in $kUpdates.Keys) {
    $ws.Cells.Item($row, 11).Value = $kUpdates[$row]
}

# --- Append newly observed board states (rows 761-804) ---
$newRows = @(
    ,@(0,0,0,0,0,0,-1,0,0,16,2)
    ,@(0,0,0,0,1,0,-1,0,-1,11,1)
    ,@(0,-1,0,0,1,0,-1,1,-1,19,1)
    ,@(1,-1,0,0,1,-1,-1,1,-1,14,1)
    ,@(-1,0,0,0,0,0,-1,0,1,15,1)
    ,@(0,0,1,0,0,0,-1,-1,0,11,2)
    ,@(0,0,1,0,0,-1,-1,-1,1,15,2)
    ,@(1,-1,-1,-1,1,0,1,1,-1,14,4)
    ,@(1,1,-1,-1,0,-1,1,0,0,11,2)
    ,@(0,0,0,0,0,0,-1,0,0,12,11)
    ,@(-1,0,0,1,0,0,-1,0,0,18,2)
    ,@(-1,-1,1,1,0,0,-1,0,0,14,2)
    ,@(-1,-1,1,1,0,1,-1,0,-1,13,2)
    ,@(1,1,-1,-1,0,0,1,0,-1,17,8)
    ,@(1,1,-1,-1,-1,0,1,1,-1,14,10)
    ,@(1,0,-1,-1,-1,1,1,1,-1,19,5)
    ,@(1,0,-1,-1,0,1,1,0,-1,19,7)
    ,@(1,1,-1,-1,0,1,1,-1,-1,13,4)
    ,@(1,1,-1,-1,-1,1,1,0,-1,17,3)
    ,@(-1,0,0,1,0,0,-1,0,0,13,2)
    ,@(-1,0,-1,1,1,0,-1,0,0,14,2)
    ,@(1,1,-1,-1,-1,0,1,-1,1,14,5)
    ,@(-1,0,-1,1,1,-1,1,-1,1,19,4)
    ,@(0,1,-1,-1,-1,1,1,1,-1,15,2)
    ,@(-1,1,-1,-1,1,1,1,-1,0,11,1)
    ,@(-1,1,-1,0,-1,1,1,-1,1,12,1)
    ,@(1,-1,-1,-1,0,1,1,1,-1,13,4)
    ,@(-1,0,0,1,0,0,-1,0,0,17,2)
    ,@(-1,0,0,1,0,0,-1,1,-1,13,2)
    ,@(-1,1,-1,1,1,-1,1,-1,0,11,3)
    ,@(-1,0,0,1,0,0,-1,0,0,11,5)
    ,@(-1,-1,0,1,0,0,-1,0,1,18,10)
    ,@(-1,-1,1,1,0,0,-1,-1,1,14,10)
    ,@(0,0,0,0,0,0,-1,0,0,11,5)
    ,@(-1,0,0,0,0,0,-1,0,1,12,5)
    ,@(1,1,-1,-1,-1,0,1,0,0,17,2)
    ,@(0,1,-1,-1,-1,1,1,-1,1,15,2)
    ,@(0,0,0,0,0,0,-1,0,0,17,1)
    ,@(0,0,-1,0,0,0,-1,1,0,13,1)
    ,@(0,-1,-1,0,1,0,-1,1,0,15,1)
    ,@(1,-1,-1,0,1,-1,-1,1,0,11,1)
    ,@(1,0,-1,-1,-1,1,1,-1,1,19,1)
    ,@(1,-1,1,-1,-1,0,1,0,0,17,1)
    ,@(1,-1,1,-1,-1,0,1,1,-1,14,1)
)

$moveText = @{
    11 = "{2, 2}"
    12 = "{1, 0}"
    13 = "{1, 1}"
    14 = "{1, 2}"
    15 = "{0, 0}"
    16 = "{2, 0}"
    17 = "{2, 1}"
    18 = "{0, 2}"
    19 = "{0, 1}"
}

$startRow = 761
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c-1]
    }
    $ws.Cells.Item($r, 10).Value = $moveText[[int]$data[9]]
    $ws.Cells.Item($r, 11).Value = $data[10]
}

"done"